$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing segment names from column A (rows 2-20) before shifting
$names = @()
for ($r = 2; $r -le 20; $r++) {
    $names += $ws.Cells.Item($r, 1).Value2
}

# Insert a new column before column B; this shifts old B,C,D -> C,D,E
$ws.Columns.Item(2).Insert()

# New header for inserted column B (copy header formatting from neighboring header cell)
$ws.Cells.Item(1, 2).Value2 = "segments"
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 2).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The data rows' new column B (holding the moved segment names) has no special
# formatting, while column A keeps its original bordered/bold style since it
# is populated with the numeric index values.
$ws.Range("B2:B20").ClearFormats()

# Fill column A with numeric index (0-based) and column B with the segment names
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 2
    $ws.Cells.Item($r, 2).Value2 = $names[$r - 2]
}
